$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# Fix row 13 date: 2019-03-21 -> 2019-03-22 (serial 43545 -> 43546)
$ws.Range("A13").Value = (Get-Date -Year 2019 -Month 3 -Day 22).Date

# Add new row 14 for problem 64 (simple dp), matching row 13's date style
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = (Get-Date -Year 2019 -Month 3 -Day 23).Date
$ws.Range("B14").Value = "64 dp"
$ws.Range("E14").Value = "done"

$ws.Range("I14").Select()
